$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Seniority (column D) values for the shuffled rows
$ws.Range("D2").Value = 1
$ws.Range("D5").Value = 5
$ws.Range("D8").Value = 3
$ws.Range("D9").Value = 3
$ws.Range("D10").Value = 7
$ws.Range("D12").Value = 29
$ws.Range("D13").Value = 9
$ws.Range("D15").Value = 65
$ws.Range("D17").Value = 6
$ws.Range("D19").Value = 2
$ws.Range("D21").Value = 5
$ws.Range("D22").Value = 3
$ws.Range("D23").Value = 65
$ws.Range("D24").Value = 9
$ws.Range("D25").Value = 7
$ws.Range("D28").Value = 6
$ws.Range("D29").Value = 0

# Move the view / selection from column C to column D
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("D29").Select()
